$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: the order in which NEW distinct strings are first written matters -
# it determines their position in the shared-string table, which must match
# the target file bit-for-bit. B1 ("New Name") and C1 ("Normalisation")
# already hold the correct text in the source workbook, so they are left
# untouched to preserve their existing shared-string slots.

# 1. Header left cell
$ws.Cells.Item(1, 1).Value = "Original Column Name"

# 2. KIV section labels (column A, rows 14-19)
$ws.Cells.Item(14, 1).Value = "KIV"
$ws.Cells.Item(15, 1).Value = "Vendor"
$ws.Cells.Item(16, 1).Value = "Tags"
$ws.Cells.Item(17, 1).Value = "Risk Level"
$ws.Cells.Item(18, 1).Value = "Source"
$ws.Cells.Item(19, 1).Value = "Billing Province Name"

# 3. Original-column names (column A, rows 2-7)
$ws.Cells.Item(2, 1).Value = "Name"
$ws.Cells.Item(3, 1).Value = "Created at"
$ws.Cells.Item(4, 1).Value = "Lineitem quantity"
$ws.Cells.Item(5, 1).Value = "Lineitem name"
$ws.Cells.Item(6, 1).Value = "Lineitem price"
$ws.Cells.Item(7, 1).Value = "Lineitem discount"

# 4. New names / normalisation notes, in the exact order they were authored
$ws.Cells.Item(2, 2).Value = "customer_id"
$ws.Cells.Item(2, 3).Value = "remove the hashtags and enumerate from 0"

$ws.Cells.Item(3, 2).Value = "order_time"

$ws.Cells.Item(6, 2).Value = "item_price"
$ws.Cells.Item(4, 2).Value = "item_quantity"

$ws.Cells.Item(3, 3).Value = "NA"

$ws.Cells.Item(6, 3).Value = "standardise from 0 to 1"
$ws.Cells.Item(5, 3).Value = "join with product data and the convert to their IDs"
$ws.Cells.Item(7, 3).Value = "change to percentage in reference to nominal item price. Then bin them to categories to interval of 20%"

$ws.Cells.Item(8, 2).Value  = "item_discount_sub_twenty"
$ws.Cells.Item(9, 2).Value  = "item_discount_sub_forty"
$ws.Cells.Item(10, 2).Value = "item_discount_sub_sixty"
$ws.Cells.Item(11, 2).Value = "item_discount_sub_eighty"
$ws.Cells.Item(12, 2).Value = "item_discount_sub_hundred"

$ws.Cells.Item(7, 2).Value = "item_discount_null"
$ws.Cells.Item(5, 2).Value = "item_id"

$ws.Cells.Item(4, 3).Value = "NA"

# --- Formatting ---
# Header row: bold + centered
$headerRange = $ws.Range("A1:C1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108  # xlCenter

# KIV title: bold only
$ws.Cells.Item(14, 1).Font.Bold = $true

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 22.21875
$ws.Columns.Item(2).ColumnWidth = 27.109375
$ws.Columns.Item(3).ColumnWidth = 22.77734375

# --- Selection ---
$ws.Range("B5").Select()
